$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Protokoll" sheet to "Protocol"
$ws.Name = "Protocol"

# Pascal Lagger's (row 5) combined progress note was split into two
# separate notes (C5 / D5), each with its own progress percentage
# tracked in the row below (C6 / D6).
$ws.Range("C5").Value = "Anpassung Datenmodell, Recherche Webtechnologien"
$ws.Range("D5").Value = "Erstellen von Triggern etc. Aufsetzen einer Testumgebung, Einpflegung des Datenmodells"
$ws.Range("D5").WrapText = $true

$ws.Range("C6").Value = 0.55
$ws.Range("D6").Value = 0

# Update selection / scroll position to match the edited view
$ws.Activate()
$ws.Range("B3").Select()
